# This sheet renders a RUML sequence diagram: columns A..G hold the vertical
# "lifeline" cells (style s=2, thick left border) plus the arrow glyphs that
# represent calls/returns, while column H holds the class/member labels.
#
# This edit inserts a new lifeline (for the newly diagrammed Vehicle/Car
# chain in rows 12-17) which shifts several existing arrow glyphs one or two
# columns to the right, and touches up a couple of arrow directions in the
# Farzi/Bike section (rows 18-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Step 1: remove cells that must disappear entirely (no style, no value)
# ---------------------------------------------------------------------
$cellsToClear = @(
    "C5", "C6", "C7",
    "C8", "C9", "C10", "C11",
    "C12", "D12",
    "C13", "D13",
    "C14", "D14",
    "C15", "D15",
    "C16", "D16",
    "C17", "D17",
    "D18",
    "D19",
    "D20",
    "D21"
)
foreach ($addr in $cellsToClear) {
    $ws.Range($addr).Clear() | Out-Null
}

# ---------------------------------------------------------------------
# Step 2: bring in the lifeline style (s=2, thick left border) for brand
# new cells that stay blank, by copying format from an existing lifeline
# cell (B2) and pasting formats only.
# ---------------------------------------------------------------------
$styleSource = $ws.Range("B2")

$emptyNewCells = @(
    "E6", "E7",
    "D9",
    "E11",
    "E12", "A13", "E13", "A14", "E14", "E15", "A16", "E16", "A17", "E17",
    "E18",
    "E19",
    "E20"
)

$styleSource.Copy() | Out-Null
foreach ($addr in $emptyNewCells) {
    $ws.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
}

# ---------------------------------------------------------------------
# Step 3: brand new lifeline cells that also carry an arrow glyph.
# ---------------------------------------------------------------------
$newValueCells = @{
    "E5"  = "⇒"
    "D8"  = "⇒"
    "D10" = "←"
    "A12" = "⇒"
    "A15" = "←"
    "E21" = "←"
}

$styleSource.Copy() | Out-Null
foreach ($addr in $newValueCells.Keys) {
    $ws.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($addr).Value = $newValueCells[$addr]
}

# ---------------------------------------------------------------------
# Step 4: existing lifeline cells whose arrow glyph is removed (kept
# blank, style untouched).
# ---------------------------------------------------------------------
$contentsToClear = @("E8", "E10", "A18")
foreach ($addr in $contentsToClear) {
    $ws.Range($addr).ClearContents() | Out-Null
}

# ---------------------------------------------------------------------
# Step 5: existing lifeline cells whose arrow glyph changes direction.
# ---------------------------------------------------------------------
$valuesToSet = @{
    "C18" = "⇒"
    "A21" = "←"
    "C21" = "◁"
}
foreach ($addr in $valuesToSet.Keys) {
    $ws.Range($addr).Value = $valuesToSet[$addr]
}
